$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $ws.Range("Z1").Formula = "=""" + $text + """"
    $ws.Range("Z1").Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

Set-TextValue "A830" "80267923"
Set-TextValue "B830" "15386-DLO-I"
$ws.Range("C830").Value = 1

Set-TextValue "A831" "80267924"
Set-TextValue "B831" "10000-LDG-I"
$ws.Range("C831").Value = 1

Set-TextValue "A832" "80267927"
Set-TextValue "B832" "33380-ATE-I"
$ws.Range("C832").Value = 10

Set-TextValue "A833" "80267927"
Set-TextValue "B833" "33381-ATE-I"
$ws.Range("C833").Value = 10

Set-TextValue "A834" "80267927"
Set-TextValue "B834" "33383-ATE-I"
$ws.Range("C834").Value = 10

Set-TextValue "A835" "80267927"
Set-TextValue "B835" "23359-ATE-I"
$ws.Range("C835").Value = 30

Set-TextValue "A836" "80267927"
Set-TextValue "B836" "10119-ATE-I"
$ws.Range("C836").Value = 3

Set-TextValue "A837" "80267927"
Set-TextValue "B837" "23364-ATE-I"
$ws.Range("C837").Value = 2

Set-TextValue "A838" "80267927"
Set-TextValue "B838" "23360-ATE-I"
$ws.Range("C838").Value = 1

Set-TextValue "A839" "80267927"
Set-TextValue "B839" "33630-ATE-I"
$ws.Range("C839").Value = 1

Set-TextValue "A840" "80267927"
Set-TextValue "B840" "33660-ATE-I"
$ws.Range("C840").Value = 1

Set-TextValue "A841" "80267927"
Set-TextValue "B841" "33378-ATE-I"
$ws.Range("C841").Value = 10

Set-TextValue "A842" "80267927"
Set-TextValue "B842" "10029-ATE-I"
$ws.Range("C842").Value = 10

Set-TextValue "A843" "80267927"
Set-TextValue "B843" "33679-ATE-I"
$ws.Range("C843").Value = 10

Set-TextValue "A844" "84004850"
Set-TextValue "B844" "10540-ARI-I"
$ws.Range("C844").Value = 4

Set-TextValue "A845" "84004855"
Set-TextValue "B845" "10065-ARI-I"
$ws.Range("C845").Value = 1

Set-TextValue "A846" "84004856"
Set-TextValue "B846" "10655-ARI-I"
$ws.Range("C846").Value = 8

Set-TextValue "A847" "84004856"
Set-TextValue "B847" "10359-ARI-I"
$ws.Range("C847").Value = 6

Set-TextValue "A848" "84004856"
Set-TextValue "B848" "10259-ARI-I"
$ws.Range("C848").Value = 4

Set-TextValue "A849" "84004856"
Set-TextValue "B849" "10257-ARI-I"
$ws.Range("C849").Value = 4

Set-TextValue "A850" "84004857"
Set-TextValue "B850" "10361-ARI-I"
$ws.Range("C850").Value = 1

Set-TextValue "A851" "84004858"
Set-TextValue "B851" "10010-ARI-I"
$ws.Range("C851").Value = 1

Set-TextValue "A852" "84004859"
Set-TextValue "B852" "10496-ARI-I"
$ws.Range("C852").Value = 3

Set-TextValue "A853" "84004859"
Set-TextValue "B853" "10496-ARI-I"
$ws.Range("C853").Value = 1

$ws.Range("Z1").Value = $null

$ws.Range("A2:C853").Select()